$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.017.87'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '3.522.23'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('D7').Value = '3.520.82'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('E11').Value = '  +3.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.386'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '4.121.90'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '3.522.47'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '65.014.74'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '392.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').Value = '3.664.50'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -3.29%  '
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.68%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('D33').Value = '3.525.01'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +6.07%  '
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '168.30'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0815'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.821'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('E43').Value = '  +5.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.96'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.21%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = '2.411.51'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('E51').Value = '  +5.56%  '
